# Updates the cryptos list prices/volumes per the Nov 7 2024 data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "75.927.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.846.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +7.09%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "193.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "598.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.58%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.552"
$ws.Range("D8").Style = "Normal"

$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.845.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.11%  "

$ws.Range("E12").Value = "  -2.08%  "

$ws.Range("E13").Value = "  +4.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.373.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.799.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.12%  "

$ws.Range("E16").Value = "  +5.02%  "

$ws.Range("E17").Value = "  +1.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.847.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.64%  "

$ws.Range("E22").Value = "  +2.80%  "

$ws.Range("E23").Value = "  +2.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.992.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000104"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.980"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "522.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.92%  "

$ws.Range("E33").Value = "  +1.22%  "

$ws.Range("E34").Value = "  +4.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.11%  "

$ws.Range("E38").Value = "  +1.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "186.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.85%  "

$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("E42").Value = "  +3.84%  "

$ws.Range("E43").Value = "  +5.63%  "

$ws.Range("E44").Value = "  +2.12%  "

$ws.Range("E45").Value = "  +4.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0887"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.30%  "

$ws.Range("E48").Value = "  +2.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.575"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.61%  "

$ws.Range("E50").Value = "  +4.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.659"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.76%  "
